$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F21").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("F29").Value = -4
$ws.Range("F32").Value = -3
$ws.Range("F34").Value = -1
